# Apply the update to the "Bill Summary" worksheet:
#  1. Insert a new row at row 14 (shifts the Totals/Grand-Total block down by one row).
#  2. Refill the new row 14 with the new BOQ line item ("Each" / switch item).
#  3. Update the quantity / amount figures that changed on the existing item rows.
#  4. Recompute the Grand Total / Net Payable figures to match the new line items.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Insert a new row above the current row 14 (pushes rows 14-20 down to 15-21) ---
$ws.Rows(14).Insert()

# --- 2. Updated quantities on the existing item rows (C column) ---
$ws.Range("C8").Value = 90
$ws.Range("C9").Value = 44
$ws.Range("C10").Value = 21
$ws.Range("C11").Value = 39
$ws.Range("C12").Value = 87
$ws.Range("C13").Value = 53

# --- Updated "Upto date Amount" figures (G column, stored as text) tied to the qty changes ---
$ws.Range("G9").NumberFormat = "@"
$ws.Range("G9").Value = "11264.00"
$ws.Range("G10").NumberFormat = "@"
$ws.Range("G10").Value = "9912.00"
$ws.Range("G11").NumberFormat = "@"
$ws.Range("G11").Value = "25818.00"
$ws.Range("G13").NumberFormat = "@"
$ws.Range("G13").Value = "7208.00"

# --- 3. Populate the newly inserted row 14 with the new BOQ line item ---
$ws.Range("A14").Value = "Each"
$ws.Range("B14").Value = 0
$ws.Range("C14").Value = 16
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "3.0"
$ws.Range("E14").Value = 'P & F ISI marked (IS:3854) 6 amp. flush type non modular switch  with CM/L no. printed and made out from industrial grade Polycarbonate or fire resistant ABS material including cutting hole in tile and making connection testing etc. as required. All as per pre approved by Engineer in charge.  For additional technical parameters of products/ work  , refer   Annexure "A" attached with this BSR .'
$ws.Range("F14").Value = 23
$ws.Range("G14").NumberFormat = "@"
$ws.Range("G14").Value = "368.00"
$ws.Range("H14").Value = 0
$ws.Range("I14").Value = ""

# --- The post-insert C column noise on the (now shifted) summary rows 15-17 ---
$ws.Range("C15").Value = 18
$ws.Range("C16").Value = 51
$ws.Range("C17").Value = 65

# --- 4. Recomputed Grand Total / Net Payable figures (rows 19 & 21 after the insert) ---
$ws.Range("G19").NumberFormat = "@"
$ws.Range("G19").Value = "54570.00"
$ws.Range("H19").NumberFormat = "@"
$ws.Range("H19").Value = "54570.00"

$ws.Range("G21").NumberFormat = "@"
$ws.Range("G21").Value = "54570.00"
$ws.Range("H21").NumberFormat = "@"
$ws.Range("H21").Value = "54570.00"
